# Apply the weekly CompStat data refresh (new crime data collected).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header strings -------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/14/2024  Through  10/20/2024"

# --- Cells that switch from a numeric value to a text placeholder ---------
# (written as text first, then re-stamped with the workbook's existing
# "placeholder" number format so no new style entries are introduced)
$styleTextSrc = $ws.Range("D15")
$ws.Range("F14").Value = "'0"
$ws.Range("F23").Value = "'0"
$styleTextSrc.Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Cells that switch from a text placeholder to a whole number ----------
$styleIntSrc = $ws.Range("G14")
$styleIntSrc.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("G30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C15").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("G29").Value = 2
$ws.Range("D30").Value = 2
$ws.Range("G30").Value = 2

# --- Cells that switch from a text placeholder to a decimal % change ------
$styleDecSrc = $ws.Range("K14")
$styleDecSrc.Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("H30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = -100
$ws.Range("E29").Value = -100
$ws.Range("H29").Value = -100
$ws.Range("E30").Value = -100
$ws.Range("H30").Value = -100

# --- Remaining numeric updates (style unchanged) ---------------------------
$ws.Range("H14").Value = -100
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("L15").Value = -15.151515151515
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 18
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 253
$ws.Range("J16").Value = 257
$ws.Range("K16").Value = -1.556420233463
$ws.Range("L16").Value = 7.659574468085
$ws.Range("M16").Value = 57.142857142857
$ws.Range("N16").Value = -74.104401228249
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -44.444444444444
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = -2.941176470588
$ws.Range("I17").Value = 418
$ws.Range("J17").Value = 362
$ws.Range("K17").Value = 15.469613259668
$ws.Range("L17").Value = 56.554307116104
$ws.Range("M17").Value = 138.857142857143
$ws.Range("N17").Value = 40.740740740740
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 13
$ws.Range("E18").Value = -46.153846153846
$ws.Range("F18").Value = 30
$ws.Range("G18").Value = 38
$ws.Range("H18").Value = -21.052631578947
$ws.Range("I18").Value = 331
$ws.Range("J18").Value = 449
$ws.Range("K18").Value = -26.280623608017
$ws.Range("L18").Value = -25.282167042889
$ws.Range("M18").Value = -16.414141414141
$ws.Range("N18").Value = -83.060388945752
$ws.Range("C19").Value = 22
$ws.Range("E19").Value = -15.384615384615
$ws.Range("F19").Value = 82
$ws.Range("G19").Value = 95
$ws.Range("H19").Value = -13.684210526315
$ws.Range("I19").Value = 1056
$ws.Range("J19").Value = 1109
$ws.Range("K19").Value = -4.779080252479
$ws.Range("L19").Value = -13.086419753086
$ws.Range("M19").Value = 99.621928166351
$ws.Range("N19").Value = -2.131603336422
$ws.Range("D20").Value = 22
$ws.Range("E20").Value = -54.545454545454
$ws.Range("F20").Value = 55
$ws.Range("G20").Value = 54
$ws.Range("H20").Value = 1.851851851851
$ws.Range("I20").Value = 450
$ws.Range("J20").Value = 416
$ws.Range("K20").Value = 8.173076923076
$ws.Range("L20").Value = 118.446601941748
$ws.Range("M20").Value = 89.873417721519
$ws.Range("N20").Value = -86.234322422759
$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 75
$ws.Range("E21").Value = -34.666666666666
$ws.Range("F21").Value = 219
$ws.Range("G21").Value = 252
$ws.Range("H21").Value = -13.095238095238
$ws.Range("I21").Value = 2540
$ws.Range("J21").Value = 2628
$ws.Range("K21").Value = -3.348554033485
$ws.Range("L21").Value = 5.701206824802
$ws.Range("M21").Value = 67.546174142480
$ws.Range("N21").Value = -66.627250032847
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -25
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = -17.857142857142
$ws.Range("L22").Value = -30.303030303030
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -100
$ws.Range("J23").Value = 14
$ws.Range("K23").Value = 35.714285714285
$ws.Range("L23").Value = 46.153846153846
$ws.Range("C24").Value = 70
$ws.Range("D24").Value = 72
$ws.Range("E24").Value = -2.777777777777
$ws.Range("F24").Value = 230
$ws.Range("G24").Value = 215
$ws.Range("H24").Value = 6.976744186046
$ws.Range("I24").Value = 2225
$ws.Range("J24").Value = 2298
$ws.Range("K24").Value = -3.176675369886
$ws.Range("L24").Value = -0.979083222073
$ws.Range("M24").Value = 77.290836653386
$ws.Range("C25").Value = 40
$ws.Range("D25").Value = 37
$ws.Range("E25").Value = 8.108108108108
$ws.Range("F25").Value = 146
$ws.Range("G25").Value = 126
$ws.Range("H25").Value = 15.873015873015
$ws.Range("I25").Value = 1469
$ws.Range("J25").Value = 1260
$ws.Range("K25").Value = 16.587301587301
$ws.Range("L25").Value = 23.549201009251
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = 17.647058823529
$ws.Range("I26").Value = 785
$ws.Range("J26").Value = 724
$ws.Range("K26").Value = 8.425414364640
$ws.Range("L26").Value = 39.679715302491
$ws.Range("M26").Value = 35.344827586206
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("L27").Value = -18.181818181818
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 78
$ws.Range("J28").Value = 84
$ws.Range("K28").Value = -7.142857142857
$ws.Range("L28").Value = 18.181818181818
$ws.Range("J29").Value = 8
$ws.Range("J30").Value = 8
$ws.Range("G31").Value = 1
$ws.Range("L31").Value = -53.846153846153
$ws.Range("L33").Value = 225
